$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B4").Value = "CommuneomVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
